$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.873.52'
$ws.Range("E2").Value = '  +1.58%  '
$ws.Range("D3").Value = '2.618.35'
$ws.Range("E3").Value = '  +1.29%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '600.99'
$ws.Range("E5").Value = '  +1.28%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '154.59'
$ws.Range("E6").Value = '  -0.20%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("E8").Value = '  +1.33%  '
$ws.Range("D9").Value = '2.616.70'
$ws.Range("E9").Value = '  +1.28%  '
$ws.Range("E10").Value = '  +10.31%  '
$ws.Range("E11").Value = '  +0.89%  '
$ws.Range("E12").Value = '  +0.86%  '
$ws.Range("E13").Value = '  -0.96%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.64'
$ws.Range("E14").Value = '  -2.13%  '
$ws.Range("E15").Value = '  +3.42%  '
$ws.Range("D16").Value = '3.097.46'
$ws.Range("E16").Value = '  +1.37%  '
$ws.Range("D17").Value = '67.737.81'
$ws.Range("E17").Value = '  +1.24%  '
$ws.Range("D18").Value = '2.613.84'
$ws.Range("E18").Value = '  +1.09%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.18'
$ws.Range("E19").Value = '  -0.87%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '366.13'
$ws.Range("E20").Value = '  +3.18%  '
$ws.Range("E21").Value = '  -1.69%  '
$ws.Range("E22").Value = '  -0.54%  '
$ws.Range("E23").Value = '  -2.20%  '
$ws.Range("E24").Value = '  -0.11%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '70.33'
$ws.Range("E25").Value = '  +4.49%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.87'
$ws.Range("E26").Value = '  -6.88%  '
$ws.Range("E27").Value = '  +0.86%  '
$ws.Range("D28").Value = '2.745.56'
$ws.Range("E28").Value = '  +1.05%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '576.82'
$ws.Range("E29").Value = '  -4.08%  '
$ws.Range("E30").Value = '  +0.40%  '
$ws.Range("E31").Value = '  -2.48%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.91'
$ws.Range("E32").Value = '  -2.02%  '
$ws.Range("E33").Value = '  +0.65%  '
$ws.Range("E35").Value = '  +0.10%  '
$ws.Range("E36").Value = '  -3.12%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.93'
$ws.Range("E37").Value = '  -1.84%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '158.21'
$ws.Range("E38").Value = '  +2.67%  '
$ws.Range("E39").Value = '  +0.76%  '
$ws.Range("E40").Value = '  +0.23%  '
$ws.Range("E41").Value = '  +3.03%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.34'
$ws.Range("E42").Value = '  -2.07%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.60'
$ws.Range("E43").Value = '  -1.96%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '41.18'
$ws.Range("E44").Value = '  -0.91%  '
$ws.Range("E45").Value = '  +0.07%  '
$ws.Range("E46").Value = '  -0.07%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '157.11'
$ws.Range("E47").Value = '  +0.84%  '
$ws.Range("D48").Value = '0.0₆0287'
$ws.Range("E48").Value = '  -7.33%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.75'
$ws.Range("E49").Value = '  +0.03%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '20.97'
$ws.Range("E50").Value = '  -1.92%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0540'
$ws.Range("E51").Value = '  -3.62%  '
